$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-CellText "D2" "37.133.28"
Set-CellText "E2" "  -0.02%  "

Set-CellText "D3" "2.049.96"
Set-CellText "E3" "  -1.46%  "

Set-CellText "E4" "  +0.06%  "

Set-CellText "D5" "247.94"
Set-CellText "E5" "  -2.20%  "

Set-CellText "E6" "  -1.60%  "

Set-CellText "D7" "57.99"
Set-CellText "E7" "  -2.26%  "

Set-CellText "E8" "  +0.04%  "

Set-CellText "E9" "  -2.61%  "

Set-CellText "D10" "0.0775"
Set-CellText "E10" "  -2.79%  "

Set-CellText "E11" "  -0.25%  "

Set-CellText "D12" "15.89"
Set-CellText "E12" "  -0.69%  "

Set-CellText "D13" "0.859"
Set-CellText "E13" "  +4.97%  "

Set-CellText "D14" "2.348.45"
Set-CellText "E14" "  -1.41%  "

Set-CellText "D15" "5.69"
Set-CellText "E15" "  +2.52%  "

Set-CellText "D16" "2.050.48"
Set-CellText "E16" "  -1.34%  "

Set-CellText "D17" "17.93"
Set-CellText "E17" "  +15.18%  "

Set-CellText "D18" "37.085.45"
Set-CellText "E18" "  -0.07%  "

Set-CellText "D19" "74.82"
Set-CellText "E19" "  +0.43%  "

Set-CellText "D20" "0.0₃0890"
Set-CellText "E20" "  -3.86%  "

Set-CellText "E21" "  -2.06%  "

Set-CellText "D22" "236.91"
Set-CellText "E22" "  -1.22%  "

Set-CellText "E23" "  -0.10%  "

Set-CellText "E24" "  +1.70%  "

Set-CellText "B25" "Monero"
Set-CellText "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText "D25" "170.01"
Set-CellText "E25" "  +0.10%  "

Set-CellText "B26" "Cosmos"
Set-CellText "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText "D26" "9.48"
Set-CellText "E26" "  +1.46%  "

Set-CellText "E27" "  -5.88%  "

Set-CellText "D28" "20.01"
Set-CellText "E28" "  -1.69%  "

Set-CellText "E29" "  -1.56%  "

Set-CellText "D30" "4.78"
Set-CellText "E30" "  -1.12%  "

Set-CellText "E31" "  -1.60%  "

Set-CellText "D32" "0.0616"
Set-CellText "E32" "  -3.23%  "

Set-CellText "E33" "  +0.59%  "

Set-CellText "E34" "  -2.10%  "

Set-CellText "E35" "  +0.07%  "

Set-CellText "E36" "  -2.40%  "

Set-CellText "D37" "1.77"
Set-CellText "E37" "  -0.34%  "

Set-CellText "D38" "3.29"
Set-CellText "E38" "  +16.80%  "

Set-CellText "E39" "  -2.74%  "

Set-CellText "E40" "  +15.37%  "

Set-CellText "D41" "0.0976"
Set-CellText "E41" "  -17.40%  "

Set-CellText "D42" "0.0222"
Set-CellText "E42" "  -2.16%  "

Set-CellText "D43" "17.26"
Set-CellText "E43" "  -3.48%  "

Set-CellText "E44" "  -2.59%  "

Set-CellText "D45" "95.63"
Set-CellText "E45" "  -3.53%  "

Set-CellText "E46" "  -0.99%  "

Set-CellText "D47" "1.274.50"
Set-CellText "E47" "  -2.53%  "

Set-CellText "E48" "  -3.93%  "

Set-CellText "E49" "  -2.16%  "

Set-CellText "D50" "2.231.23"
Set-CellText "E50" "  -1.45%  "

Set-CellText "D51" "43.79"
Set-CellText "E51" "  -0.38%  "
